$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '30.691.66'
$ws.Range('E2').NumberFormat = '@'
$ws.Range('E2').Value = '  +1.68%  '

$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.899.32'
$ws.Range('E3').NumberFormat = '@'
$ws.Range('E3').Value = '  +2.30%  '

$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.0000'
$ws.Range('E4').NumberFormat = '@'
$ws.Range('E4').Value = '  -0.25%  '

$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '239.18'
$ws.Range('E5').NumberFormat = '@'
$ws.Range('E5').Value = '  +1.20%  '

$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.9996'
$ws.Range('E6').NumberFormat = '@'
$ws.Range('E6').Value = '  -0.29%  '

$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '0.4797'
$ws.Range('E7').NumberFormat = '@'
$ws.Range('E7').Value = '  +0.59%  '

$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.2833'
$ws.Range('E8').NumberFormat = '@'
$ws.Range('E8').Value = '  +0.39%  '

$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.06541'
$ws.Range('E9').NumberFormat = '@'
$ws.Range('E9').Value = '  +0.36%  '

$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '1.933.46'
$ws.Range('E10').NumberFormat = '@'
$ws.Range('E10').Value = '  +4.05%  '

$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.07474'
$ws.Range('E11').NumberFormat = '@'
$ws.Range('E11').Value = '  +1.69%  '

$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '16.66'
$ws.Range('E12').NumberFormat = '@'
$ws.Range('E12').Value = '  +1.54%  '

$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '5.093'
$ws.Range('E13').NumberFormat = '@'
$ws.Range('E13').Value = '  -1.17%  '

$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '87.98'
$ws.Range('E14').NumberFormat = '@'
$ws.Range('E14').Value = '  +0.76%  '

$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.6676'
$ws.Range('E15').NumberFormat = '@'
$ws.Range('E15').Value = '  +3.25%  '

$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '30.663.20'
$ws.Range('E16').NumberFormat = '@'
$ws.Range('E16').Value = '  +1.72%  '

$ws.Range('E17').NumberFormat = '@'
$ws.Range('E17').Value = '  +0.47%  '

$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '0.9995'
$ws.Range('E18').NumberFormat = '@'
$ws.Range('E18').Value = '  -0.50%  '

$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '2.200.54'
$ws.Range('E19').NumberFormat = '@'
$ws.Range('E19').Value = '  +3.93%  '

$ws.Range('E20').NumberFormat = '@'
$ws.Range('E20').Value = '  -0.31%  '

$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '228.82'
$ws.Range('E21').NumberFormat = '@'
$ws.Range('E21').Value = '  +4.95%  '

$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '5.308'
$ws.Range('E22').NumberFormat = '@'
$ws.Range('E22').Value = '  +0.50%  '

$ws.Range('E23').NumberFormat = '@'
$ws.Range('E23').Value = '  -0.06%  '

$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '6.220'
$ws.Range('E24').NumberFormat = '@'
$ws.Range('E24').Value = '  +1.66%  '

$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '169.18'
$ws.Range('E25').NumberFormat = '@'
$ws.Range('E25').Value = '  +1.83%  '

$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '9.267'
$ws.Range('E26').NumberFormat = '@'
$ws.Range('E26').Value = '  -0.52%  '

$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.58'
$ws.Range('E27').NumberFormat = '@'
$ws.Range('E27').Value = '  +0.31%  '

$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '1.953'
$ws.Range('E28').NumberFormat = '@'
$ws.Range('E28').Value = '  +2.10%  '

$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.403'
$ws.Range('E29').NumberFormat = '@'
$ws.Range('E29').Value = '  -1.63%  '

$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '0.09875'
$ws.Range('E30').NumberFormat = '@'
$ws.Range('E30').Value = '  +8.26%  '

$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.345'
$ws.Range('E31').NumberFormat = '@'
$ws.Range('E31').Value = '  +1.99%  '

$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '4.019'
$ws.Range('E32').NumberFormat = '@'
$ws.Range('E32').Value = '  +1.03%  '

$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05055'
$ws.Range('E33').NumberFormat = '@'
$ws.Range('E33').Value = '  +0.23%  '

$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.224'
$ws.Range('E34').NumberFormat = '@'
$ws.Range('E34').Value = '  +7.64%  '

$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '0.7523'
$ws.Range('E35').NumberFormat = '@'
$ws.Range('E35').Value = '  +1.02%  '

$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '2.711'
$ws.Range('E36').NumberFormat = '@'
$ws.Range('E36').Value = '  +0.73%  '

$ws.Range('D37').NumberFormat = '@'
$ws.Range('D37').Value = '0.01871'
$ws.Range('E37').NumberFormat = '@'
$ws.Range('E37').Value = '  +2.54%  '

$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.654'
$ws.Range('E38').NumberFormat = '@'
$ws.Range('E38').Value = '  +1.70%  '

$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '0.9203'
$ws.Range('E39').NumberFormat = '@'
$ws.Range('E39').Value = '  +1.47%  '

$ws.Range('E40').NumberFormat = '@'
$ws.Range('E40').Value = '  +1.48%  '

$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '107.26'
$ws.Range('E41').NumberFormat = '@'
$ws.Range('E41').Value = '  +0.11%  '

$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '5.835'
$ws.Range('E42').NumberFormat = '@'
$ws.Range('E42').Value = '  -1.31%  '

$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.4290'
$ws.Range('E43').NumberFormat = '@'
$ws.Range('E43').Value = '  +0.72%  '

$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '1.006'
$ws.Range('E44').NumberFormat = '@'
$ws.Range('E44').Value = '  +0.42%  '

$ws.Range('D45').NumberFormat = '@'
$ws.Range('D45').Value = '7.395'
$ws.Range('E45').NumberFormat = '@'
$ws.Range('E45').Value = '  -0.69%  '

$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '64.52'
$ws.Range('E46').NumberFormat = '@'
$ws.Range('E46').Value = '  +0.33%  '

$ws.Range('E47').NumberFormat = '@'
$ws.Range('E47').Value = '  -3.20%  '

$ws.Range('B48').NumberFormat = '@'
$ws.Range('B48').Value = 'NEARProtocol'
$ws.Range('C48').NumberFormat = '@'
$ws.Range('C48').Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '1.481'
$ws.Range('E48').NumberFormat = '@'
$ws.Range('E48').Value = '  -5.49%  '

$ws.Range('B49').NumberFormat = '@'
$ws.Range('B49').Value = 'EnergySwap'
$ws.Range('C49').NumberFormat = '@'
$ws.Range('C49').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '8.958'
$ws.Range('E49').NumberFormat = '@'
$ws.Range('E49').Value = '  +1.41%  '

$ws.Range('E50').NumberFormat = '@'
$ws.Range('E50').Value = '  -1.39%  '

$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.05662'
$ws.Range('E51').NumberFormat = '@'
$ws.Range('E51').Value = '  -0.68%  '
